{"js": "// Replace each two-digit-multiplication answer string in the document's\n// table cells with its updated value, per the authoring diff.\n// Each old value is unique within the document, so a direct search +\n// whole-text replace for each pair is safe and order independent.\n\nconst replacements = [\n  [\"26\u00d735=910\", \"26\u00d761=1586\"],\n  [\"44\u00d791=4004\", \"25\u00d788=2200\"],\n  [\"17\u00d721=357\", \"35\u00d742=1470\"],\n  [\"92\u00d797=8924\", \"32\u00d713=416\"],\n  [\"55\u00d719=1045\", \"30\u00d764=1920\"],\n  [\"77\u00d728=2156\", \"89\u00d764=5696\"],\n  [\"99\u00d770=6930\", \"51\u00d760=3060\"],\n  [\"78\u00d794=7332\", \"77\u00d712=924\"],\n  [\"77\u00d757=4389\", \"38\u00d734=1292\"],\n  [\"39\u00d742=1638\", \"41\u00d722=902\"],\n  [\"95\u00d725=2375\", \"47\u00d778=3666\"],\n  [\"92\u00d795=8740\", \"84\u00d749=4116\"],\n  [\"44\u00d725=1100\", \"37\u00d789=3293\"],\n  [\"90\u00d730=2700\", \"82\u00d754=4428\"],\n  [\"94\u00d726=2444\", \"81\u00d772=5832\"],\n  [\"49\u00d717=833\", \"55\u00d749=2695\"],\n  [\"74\u00d787=6438\", \"14\u00d721=294\"],\n  [\"83\u00d750=4150\", \"80\u00d727=2160\"],\n  [\"16\u00d785=1360\", \"72\u00d751=3672\"],\n  [\"86\u00d766=5676\", \"27\u00d759=1593\"],\n  [\"47\u00d788=4136\", \"21\u00d724=504\"],\n  [\"97\u00d782=7954\", \"28\u00d792=2576\"],\n  [\"17\u00d773=1241\", \"83\u00d748=3984\"],\n  [\"33\u00d778=2574\", \"44\u00d751=2244\"],\n  [\"45\u00d711=495\", \"15\u00d741=615\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication answer string in the document's\n# table cells with its updated value, per the authoring diff.\n# Each old value is unique within the document, so a Find/Replace over\n# the whole document content for each pair is safe and order independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"26\u00d735=910\", \"26\u00d761=1586\"),\n  @(\"44\u00d791=4004\", \"25\u00d788=2200\"),\n  @(\"17\u00d721=357\", \"35\u00d742=1470\"),\n  @(\"92\u00d797=8924\", \"32\u00d713=416\"),\n  @(\"55\u00d719=1045\", \"30\u00d764=1920\"),\n  @(\"77\u00d728=2156\", \"89\u00d764=5696\"),\n  @(\"99\u00d770=6930\", \"51\u00d760=3060\"),\n  @(\"78\u00d794=7332\", \"77\u00d712=924\"),\n  @(\"77\u00d757=4389\", \"38\u00d734=1292\"),\n  @(\"39\u00d742=1638\", \"41\u00d722=902\"),\n  @(\"95\u00d725=2375\", \"47\u00d778=3666\"),\n  @(\"92\u00d795=8740\", \"84\u00d749=4116\"),\n  @(\"44\u00d725=1100\", \"37\u00d789=3293\"),\n  @(\"90\u00d730=2700\", \"82\u00d754=4428\"),\n  @(\"94\u00d726=2444\", \"81\u00d772=5832\"),\n  @(\"49\u00d717=833\", \"55\u00d749=2695\"),\n  @(\"74\u00d787=6438\", \"14\u00d721=294\"),\n  @(\"83\u00d750=4150\", \"80\u00d727=2160\"),\n  @(\"16\u00d785=1360\", \"72\u00d751=3672\"),\n  @(\"86\u00d766=5676\", \"27\u00d759=1593\"),\n  @(\"47\u00d788=4136\", \"21\u00d724=504\"),\n  @(\"97\u00d782=7954\", \"28\u00d792=2576\"),\n  @(\"17\u00d773=1241\", \"83\u00d748=3984\"),\n  @(\"33\u00d778=2574\", \"44\u00d751=2244\"),\n  @(\"45\u00d711=495\", \"15\u00d741=615\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
